$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values as text (e.g. "240.30", "26.721.78").
# Force Text number format on D2:D51 first so Excel does not
# auto-convert numeric-looking strings into actual numbers and
# lose formatting such as trailing zeros (matches inlineStr cells).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.721.78'
$ws.Range("E2").Value = '  +1.32%  '
$ws.Range("D3").Value = '1.725.97'
$ws.Range("E3").Value = '  +0.31%  '
$ws.Range("D4").Value = '0.9981'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '240.30'
$ws.Range("E5").Value = '  -0.63%  '
$ws.Range("D6").Value = '0.9989'
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("D7").Value = '0.4804'
$ws.Range("E7").Value = '  -1.50%  '
$ws.Range("D8").Value = '0.2585'
$ws.Range("E8").Value = '  -0.27%  '
$ws.Range("D9").Value = '0.06174'
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("D10").Value = '1.723.83'
$ws.Range("E10").Value = '  +0.18%  '
$ws.Range("D11").Value = '15.79'
$ws.Range("E11").Value = '  +1.94%  '
$ws.Range("D12").Value = '0.06864'
$ws.Range("E12").Value = '  -1.52%  '
$ws.Range("D13").Value = '0.6019'
$ws.Range("E13").Value = '  +0.94%  '
$ws.Range("D14").Value = '4.456'
$ws.Range("E14").Value = '  -1.33%  '
$ws.Range("D15").Value = '76.84'
$ws.Range("E15").Value = '  -0.25%  '
$ws.Range("D16").Value = '0.9986'
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").Value = '26.548.75'
$ws.Range("E17").Value = '  +0.67%  '
$ws.Range("D18").Value = '0.9982'
$ws.Range("E18").Value = '  -0.13%  '
$ws.Range("D19").Value = '0.000007134'
$ws.Range("E19").Value = '  -0.48%  '
$ws.Range("D20").Value = '11.35'
$ws.Range("E20").Value = '  +0.42%  '
$ws.Range("D21").Value = '1.946.31'
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("D22").Value = '4.413'
$ws.Range("E22").Value = '  -0.53%  '
$ws.Range("D23").Value = '8.510'
$ws.Range("E23").Value = '  +0.41%  '
$ws.Range("D24").Value = '5.048'
$ws.Range("E24").Value = '  -0.38%  '
$ws.Range("D25").Value = '139.75'
$ws.Range("E25").Value = '  +1.44%  '
$ws.Range("D26").Value = '15.20'
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("D27").Value = '1.769'
$ws.Range("E27").Value = '  +2.57%  '
$ws.Range("D28").Value = '106.07'
$ws.Range("E28").Value = '  -0.13%  '
$ws.Range("D29").Value = '1.366'
$ws.Range("E29").Value = '  -2.50%  '
$ws.Range("D30").Value = '3.999'
$ws.Range("E30").Value = '  +2.64%  '
$ws.Range("D31").Value = '0.07901'
$ws.Range("E31").Value = '  -1.24%  '
$ws.Range("D32").Value = '3.654'
$ws.Range("E32").Value = '  +0.15%  '
$ws.Range("D33").Value = '0.04520'
$ws.Range("E33").Value = '  +0.58%  '
$ws.Range("D34").Value = '2.597'
$ws.Range("E34").Value = '  -0.40%  '
$ws.Range("D35").Value = '0.9968'
$ws.Range("E35").Value = '  +0.19%  '
$ws.Range("D36").Value = '0.6161'
$ws.Range("E36").Value = '  -0.76%  '
$ws.Range("D37").Value = '0.9303'
$ws.Range("E37").Value = '  +0.99%  '
$ws.Range("D38").Value = '2.450'
$ws.Range("E38").Value = '  +2.79%  '
$ws.Range("D39").Value = '1.987'
$ws.Range("E39").Value = '  +1.39%  '
$ws.Range("D40").Value = '0.9980'
$ws.Range("E40").Value = '  -0.09%  '
$ws.Range("D41").Value = '0.01495'
$ws.Range("E41").Value = '  +1.40%  '
$ws.Range("D42").Value = '5.606'
$ws.Range("E42").Value = '  +3.56%  '
$ws.Range("D43").Value = '99.90'
$ws.Range("E43").Value = '  -0.32%  '
$ws.Range("D44").Value = '0.3818'
$ws.Range("E44").Value = '  -0.42%  '
$ws.Range("D45").Value = '6.745'
$ws.Range("E45").Value = '  -2.28%  '
$ws.Range("D46").Value = '0.1153'
$ws.Range("E46").Value = '  -0.75%  '
$ws.Range("D47").Value = '0.05357'
$ws.Range("E47").Value = '  -0.14%  '
$ws.Range("D48").Value = '7.886'
$ws.Range("E48").Value = '  +3.01%  '
$ws.Range("D49").Value = '30.04'
$ws.Range("E49").Value = '  -0.42%  '
$ws.Range("D50").Value = '1.245'
$ws.Range("E50").Value = '  +2.24%  '
$ws.Range("D51").Value = '51.33'
$ws.Range("E51").Value = '  +0.79%  '
